$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper to write a full data record (matching the structure used throughout
# this worksheet) into a given row number. Columns A,B,C,E,F,G,H,I,N,O,Q,R are
# constant across the whole dataset; only D (Fecha), J (Volumen), K (Precio
# minimo), M (Precio promedio ponderado) and P (Precio $/Kg) vary. L (Precio
# maximo) is always 500.
function Set-EspinacaRow($row, $fecha, $volumen, $precioMin, $precioProm, $precioKg) {
    $ws.Cells.Item($row, 1).Value2 = 8
    $ws.Cells.Item($row, 2).Value2 = "Terminal La Palmera de La Serena"
    $ws.Cells.Item($row, 3).Value2 = "Coquimbo"
    $ws.Cells.Item($row, 4).Value2 = $fecha
    $ws.Cells.Item($row, 5).Value2 = 4
    $ws.Cells.Item($row, 6).Value2 = 100112012
    $ws.Cells.Item($row, 7).Value2 = "Espinaca"
    $ws.Cells.Item($row, 8).Value2 = "Sin especificar"
    $ws.Cells.Item($row, 9).Value2 = "Primera"
    $ws.Cells.Item($row, 10).Value2 = $volumen
    $ws.Cells.Item($row, 11).Value2 = $precioMin
    $ws.Cells.Item($row, 12).Value2 = 500
    $ws.Cells.Item($row, 13).Value2 = $precioProm
    $ws.Cells.Item($row, 14).Value2 = "`$/atado 300 a 500 gramos"
    $ws.Cells.Item($row, 15).Value2 = "Provincia del Elqu" + [char]0x00ED
    $ws.Cells.Item($row, 16).Value2 = $precioKg
    $ws.Cells.Item($row, 17).Value2 = 0.5
    $ws.Cells.Item($row, 18).Value2 = "Hortaliza"
}

# Insert a new weekly record before the current first row of the block (row 127),
# pushing all the existing rows (127-206) down by one.
$ws.Rows.Item(127).Insert()
Set-EspinacaRow 127 44567 3000 400 450 900

# Insert a second new weekly record further down. After the first insert, the
# target location is now row 200 (this was row 199 prior to the first insert).
$ws.Rows.Item(200).Insert()
Set-EspinacaRow 200 44568 3000 400 450 900
